$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Accountant- Family Office'
$ws.Range("B2").Value = 'The Quest Organization'
$ws.Range("C2").Value = 'New York, NY'
$ws.Range("D2").Value = 'Full-Time'
$ws.Range("E2").Value = 'Strong knowledge of ground-up stick frame (podium, wrap, garden) multifamily projects * 6+ years construction experience as a Construction Manager or Superintendent * Experience with managing ...'
$ws.Range("F2").Value = 'https://www.ziprecruiter.com/k/l/AAKn0EU9NZ9JGogmToGLK8Jm1iRzYvpzOQN41N1hst6Afychm3IZdKY25Pt-1ZxoPZLzEOgT8MzcBfbNhonWOneoHNYM1DlDU7PisiJcZFz9yQlabUeI1eRlgqIYbN4L6UgRL27BUFe7qWK8pt9B5TvZ5TOzP6yuO40CW2_e3oJalPAeQ_na'

$ws.Range("A3").Value = 'Accountant'
$ws.Range("B3").Value = 'National Systems America LP'
$ws.Range("C3").Value = 'Manhattan, NY'
$ws.Range("D3").Value = 'Full-Time'
$ws.Range("E3").Value = 'We are seeking a Commercial Construction Superintendent to join our team! You will oversee project planning, scheduling, subcontractors work. Responsibilities: * Oversee all aspects of construction ...'
$ws.Range("F3").Value = 'https://www.ziprecruiter.com/k/l/AAIvuvup4F97ByPr-Drj6tOWUr_3mr6YwW6F-w-JonEuZfK7KhQinV_SLxolxkuYgmO49fGs7aTMbQ2u-thRQNM_7TUG6m9ZXGHxK_Fw28oDsxn6UsdB4KdRgUgQU-9MQdISvpIAypVF1S_wjFyOJVoMa06Wmon-JLDo7qjKSa1yd61nlT0B'

$ws.Range("A4").Value = 'Property Accountant'
$ws.Range("B4").Value = 'Creative Financial Staffing (CFS)'
$ws.Range("C4").Value = 'Manhattan, NY'
$ws.Range("D4").Value = 'Full-Time'
$ws.Range("E4").Value = 'Education * 4 yr degree in engineering, construction science, construction management or similar Position Summary We are looking for an Assistant Project Manager to perform the technical and support ...'
$ws.Range("F4").Value = 'https://www.ziprecruiter.com/k/l/AAJuEZjmoIXW-k7Dlm0LiYTxuonnmxZwGAqa-Z3dJXsUAV13uQiEnLi9f2UZvMLabNugiIZTJkCJ1waEZlap7kKPpVBNaaTY6mKm-wORFafli8q6jxs4Tqgo7DzQHyNA6_QYjP2MN1gJ87mFvOlAoSOPas4p3kWGAf8LvDG2F_8YWLwPPHcY'

$ws.Range("A5").Value = 'Bookkeeper Accountant'
$ws.Range("B5").Value = 'Abetta Boiler Welding INC'
$ws.Range("C5").Value = 'Manhattan, NY'
$ws.Range("D5").Value = 'Full-Time'
$ws.Range("E5").Value = 'Talent Corps is currently hiring for General Labors for the Fort Worth, TX area for a Temp to Perm position. Candidates are needing to be physically fit as this job is physically demanding and may ...'
$ws.Range("F5").Value = 'https://www.ziprecruiter.com/k/l/AAK-NMF07hninHmQnqOLp9jlHULdi_6fGcx8gSDSyqvcwf9BpzW-w9HmASQzm22TlmkXYgFvuhJusRoEE55hW_3G1NgiwcTOCj_H9THhoLznCiAGLbcC3yO6peQYL0EeeJVgiHawOdexw7wvp92NPmlG7GzaW5nlSE_o1p-oQ_IgdHAfDwKS'

$ws.Range("A6").Value = 'Staff Accountant'
$ws.Range("B6").Value = 'Indcorp Fiscal Services Inc.'
$ws.Range("C6").Value = 'Island Park, NY'
$ws.Range("D6").Value = 'Full-Time'
$ws.Range("E6").Value = 'Restaurant Depot is seeking a Construction & Maintenance Manager. Position will be based out of the Dallas/Richardson area offices. This person will be working under the direction of the Director of ...'
$ws.Range("F6").Value = 'https://www.ziprecruiter.com/k/l/AAL-AMNGLKjtIy0ngvz2m-3vpBPd1s52zBx6-VE6JLDeOKKAPVn7I1ewhpaIRHdWKQpAk477uSje-L9QKsm-C15YmZ0_3gp6Y5Nuy1X0rgnBimX6npM4FrFjrJo7viKhMLMAXmRip9nVZNFmU0wp6MFmtqb0I40POgJ_RzELGjcda2AALpMW'

$ws.Range("A7").Value = 'Junior Accountant'
$ws.Range("B7").Value = 'Greater New York Insurance Companies'
$ws.Range("C7").Value = 'Manhattan, NY'
$ws.Range("D7").Value = 'Full-Time'
$ws.Range("E7").Value = 'Texas Home Builder seeking Construction Managers This Jobot Job is hosted by: Henry Chan Are you a fit? Easy Apply now by clicking the "Apply Now" button and sending us your resume. Salary: $65,000 ...'
$ws.Range("F7").Value = 'https://www.ziprecruiter.com/k/l/AAIM8p5PnmW1w3qLf-pRhZP-tax17jRIloz7G5phthY7wi-StrjG0yb0Di35N-URKP4Naq8cirV5VRYfW9w7WFjt471WZs7qRKNv72y0mhUSBa5ZTG_kGAMf720CRM9APtmLvgBCYv8KbSaO-OWYC2ffX0cKFcH2HXPwgTKHM5OkOX5qlbgi'

$ws.Range("A8").Value = 'Accountant'
$ws.Range("B8").Value = 'E-Z TAX GROUP NORTH AMERICA INC'
$ws.Range("C8").Value = 'Flushing, NY'
$ws.Range("D8").Value = 'Full-Time'
$ws.Range("E8").Value = 'The Construction Manager role is a player/coach position responsible for both managing overall construction execution of individual projects and overseeing 3rd party construction management firm in ...'
$ws.Range("F8").Value = 'https://www.ziprecruiter.com/k/l/AAKcAIY5P2hmJb5reyjOKUvhwskuRvOj5yVtB5HQiOaO4T37Lxr_jEngwbYJq4btvK5Wg_LeKZIjF76gvcJa6o66L2CXLQQwfpzc0FKHpzTd1kPB84_X242b5gwETQlpIFyl7F-jGBG-pDko7DKMRaggkwxwSkxYYibXmRPay6q0LpnUApgU'

$ws.Range("A9").Value = 'Staff Accountant'
$ws.Range("B9").Value = 'Real Estate Management Company'
$ws.Range("C9").Value = 'Bellmore, NY'
$ws.Range("D9").Value = 'Full-Time'
$ws.Range("E9").Value = 'Primary Function The General Superintendent oversees the construction of multiple projects. He/she will coordinate with other superintendents to ensure other jobsites are being manned properly and ...'
$ws.Range("F9").Value = 'https://www.ziprecruiter.com/k/l/AALHIk41TDxgz9HrNt5nNVuCqj5RpqW32nBI-4ahi76pFtyVNldXPXtgef9EsSvmbt6CTkQcCBBLgEYeXze2EvKZfIuhPdq9u3kDC6mEQYS56K2PdX_QbMX3lHEQs_GEdgAugfk_Zj2e5h7HcZr1wftfXZKgUYanc67axITWgHnBRHB6MxTM'

$ws.Range("A10").Value = 'Staff Accountant'
$ws.Range("B10").Value = 'Zing Professional Recruiting a Peoplelink Group Company'
$ws.Range("C10").Value = 'Norwood, NJ'
$ws.Range("D10").Value = 'Full-Time'
$ws.Range("E10").Value = 'We are seeking a Construction Project Engineer to join our team! Responsibilities: * Verify project drawing and specification compliance * Method of Procedure research and development * Construction ...'
$ws.Range("F10").Value = 'https://www.ziprecruiter.com/k/l/AAKeJAWObIEBRL7z0gUgVHe5ao4Y3_vwoloBtiAXsnPjIRCpqeE5XxA0tzkmVWmnbkHu0e8BEd4xe2g-i_7EMWWpEKWhZiCoxt_VM1T0ciLOak40lIEJAm773F9logEkzmtJSOqbemFFVcrvW6HOxysvAa1flIzMjE6m6rkj39xanK_FnD5U'

$ws.Range("A11").Value = 'Senior Accountant, Ecommerce Clients'
$ws.Range("B11").Value = 'Accountingfly'
$ws.Range("C11").Value = 'New York, New York'
$ws.Range("D11").Value = 'Full-Time'
$ws.Range("E11").Value = 'Perform drywall, finishing, and Painting to complete tenant improvements for national retail construction. Truck and tools are required. Extensive travel is required with this position. Company ...'
$ws.Range("F11").Value = 'https://www.ziprecruiter.com/ek/l/AAKjdTbKJGppBgL4HYKoTjVnMU-iYj76Lm_6l344SXV4LnKbwFQ3jucpGrEZLpkRfjJRGiLXUmQSysus-MNW4ZOiGl9jwzw_LDW6-tyJgLAPSmsPzdyWnQwb6-3W4o_tvcULFzW9APOUx3UiTVL_uvYUeQv0L47Yv4ecDg7O8DdbMTFg6vQgmg'
